# Generate Report for Handback
# Applies the "handed back" status update to the localization-status workbook:
#   - flips the "In Translation" status to "Handed back: in sync with en-US"
#   - fills in the Latest Target File / Latest Handback File / Latest Handback
#     DateTime columns for each language sheet, with hyperlinks on the target
#     file cells
#   - widens a handful of columns so the new long values are readable

$wb = $excel.ActiveWorkbook

$oldStatus = "In Translation"
$newStatus = "Handed back: in sync with en-US"

$b4107e81Url = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/d849b26eead8f8f0675dd25fedd19cf1183acb40/e2e/b4107e81-858c-4b4b-9f0b-2c453916d44e.md"
$dbe0ca11Url = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/d849b26eead8f8f0675dd25fedd19cf1183acb40/e2e/dbe0ca11-ebea-4c0c-b8f3-c955ec718bbd.md"

# ---------------------------------------------------------------------------
# Overview sheet: flip the per-language status cells and widen the columns
# ---------------------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")

foreach ($addr in @("E2", "F2", "E3", "F3")) {
    if ($overview.Range($addr).Value2 -eq $oldStatus) {
        $overview.Range($addr).Value = $newStatus
    }
}

$overview.Columns.Item(5).ColumnWidth = 29.9777047293527
$overview.Columns.Item(6).ColumnWidth = 29.9777047293527

# ---------------------------------------------------------------------------
# Per-language detail sheets (zh-cn / de-de)
# ---------------------------------------------------------------------------
function Update-LanguageSheet($SheetName, $HandbackDateTime) {

    $ws = $wb.Worksheets.Item($SheetName)

    foreach ($addr in @("C2", "C3")) {
        if ($ws.Range($addr).Value2 -eq $oldStatus) {
            $ws.Range($addr).Value = $newStatus
        }
    }

    # Row 2 -> b4107e81-858c-4b4b-9f0b-2c453916d44e
    $ws.Range("I2").Value = "b4107e81-858c-4b4b-9f0b-2c453916d44e.md"
    $ws.Hyperlinks.Add($ws.Range("I2"), $b4107e81Url, [Type]::Missing, [Type]::Missing, "b4107e81-858c-4b4b-9f0b-2c453916d44e.md") | Out-Null
    $ws.Range("I2").Font.Underline = 2
    $ws.Range("I2").Font.Color = 15570276
    $ws.Range("J2").Value = "b4107e81-858c-4b4b-9f0b-2c453916d44e.eb169921748a0663e1ee62c2b8b99af4db7bd76d." + $SheetName + ".xlf"
    $ws.Range("K2").Value = $HandbackDateTime

    # Row 3 -> dbe0ca11-ebea-4c0c-b8f3-c955ec718bbd
    $ws.Range("I3").Value = "dbe0ca11-ebea-4c0c-b8f3-c955ec718bbd.md"
    $ws.Hyperlinks.Add($ws.Range("I3"), $dbe0ca11Url, [Type]::Missing, [Type]::Missing, "dbe0ca11-ebea-4c0c-b8f3-c955ec718bbd.md") | Out-Null
    $ws.Range("I3").Font.Underline = 2
    $ws.Range("I3").Font.Color = 15570276
    $ws.Range("J3").Value = "dbe0ca11-ebea-4c0c-b8f3-c955ec718bbd.e740855797a25a8b464ed71a9de82db468b77d39." + $SheetName + ".xlf"
    $ws.Range("K3").Value = $HandbackDateTime

    $ws.Columns.Item(3).ColumnWidth = 29.9777047293527
    $ws.Columns.Item(9).ColumnWidth = 40
    $ws.Columns.Item(10).ColumnWidth = 40
}

Update-LanguageSheet "zh-cn" "2016-08-22 02:36:26"
Update-LanguageSheet "de-de" "2016-08-22 02:36:32"
